{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies two template-placeholder updates to the SKU document:\n//  1. \"{islam}\" -> \"{agama}\" (and drops the spell-check proofErr markup\n//     that wrapped the old \"islam\" run).\n//  2. The hard-coded date run \", 04 JULI  2025\" is split into a literal\n//     \", \" run followed by a new \"{tanggal}\" placeholder run.\n\nconst body = context.document.body;\n\n// --- Change 1: {islam} -> {agama} -------------------------------------\nconst religionResults = body.search(\"{islam}\", { matchCase: true });\nreligionResults.load(\"text\");\nawait context.sync();\n\nif (religionResults.items.length > 0) {\n  const religionRange = religionResults.items[0];\n  const religionOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>agama</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  religionRange.insertOoxml(religionOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2: \", 04 JULI  2025\" -> \", \" + \"{tanggal}\" -----------------\nconst dateResults = body.search(\", 04 JULI  2025\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  const dateRange = dateResults.items[0];\n  const dateOoxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>{tanggal}</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  dateRange.insertOoxml(dateOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies two template-placeholder updates to the SKU document:\n#  1. \"{islam}\" -> \"{agama}\" (and drops the spell-check proofErr markup\n#     that wrapped the old \"islam\" run).\n#  2. The hard-coded date run \", 04 JULI  2025\" is split into a literal\n#     \", \" run followed by a new \"{tanggal}\" placeholder run.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: {islam} -> {agama} --------------------------------------\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\"islam\")\nif ($found1) {\n  $xml1 = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p w:rsidR=\"00B973EE\" w:rsidRPr=\"00B973EE\" w:rsidRDefault=\"004471E3\" w:rsidP=\"00B973EE\">' +\n    '<w:pPr><w:spacing w:line=\"360\" w:lineRule=\"auto\"/><w:jc w:val=\"both\"/></w:pPr>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>agama</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  $range1.InsertXML($xml1)\n}\n\n# --- Change 2: \", 04 JULI  2025\" -> \", \" + \"{tanggal}\" ------------------\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\", 04 JULI  2025\")\nif ($found2) {\n  $xml2 = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p w:rsidR=\"009451E0\" w:rsidRPr=\"00B973EE\" w:rsidRDefault=\"009451E0\" w:rsidP=\"00F37F5C\">' +\n    '<w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r w:rsidRPr=\"00B973EE\"><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Limo Koto</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>{tanggal}</w:t></w:r>' +\n    '</w:p></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  $range2.InsertXML($xml2)\n}\n"}
